$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Non Cash Items (Other)
$ws.Range("B4").Value = 341000000.0
$ws.Range("F4").Value = 181000000.0

# Row 6: Change in inventories
$ws.Range("B6").Value = -125000000.0
$ws.Range("C6").Value = -158000000.0
$ws.Range("D6").Value = -346000000.0
$ws.Range("E6").Value = -262000000.0
$ws.Range("F6").Value = -107000000.0
$ws.Range("G6").Value = -106000000.0

# Row 8: Change in payables and accrued liability
$ws.Range("B8").Value = 2324000000.0
$ws.Range("C8").Value = 3248000000.0
$ws.Range("D8").Value = 3188000000.0
$ws.Range("E8").Value = 2402000000.0
$ws.Range("F8").Value = 1652000000.0
$ws.Range("G8").Value = 582000000.0

# Row 10: Change in other assets and liabilities
$ws.Range("B10").Value = 98000000.0

# Row 32: Capital Stock Change (was an empty inline string, now a number)
$ws.Range("B32").Value = -617000000.0
